$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = -1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = -4
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -1
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = 3
$ws.Range("F18").Value = 7
$ws.Range("F19").Value = 6
$ws.Range("F20").Value = 2
$ws.Range("F21").Value = 4
$ws.Range("F22").Value = -7
$ws.Range("F23").Value = -3
$ws.Range("F24").Value = -2
$ws.Range("F25").Value = -2
$ws.Range("F26").Value = -1
$ws.Range("F28").Value = -3
$ws.Range("F29").Value = -2
